$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for B..BO (66 data columns), matching header row order
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG","BH","BI","BJ","BK","BL","BM","BN","BO")

# New "Base" row (row 2) values - replace existing halved values
$row2Values = @(13511034093.70666,13436107206.3887,13385419306.00335,12316684873.22743,11430871479.97902,11385231230.8767,11338294179.50647,11292844962.14668,11246104781.76936,11304395507.51397,11266007415.25349,11229476015.5672,11195206102.84851,11159940051.22204,11125770037.37515,11090606938.2172,11056016451.04463,11021996285.58852,10986987615.28494,10953067348.67183,10918161631.01235,11004070467.28123,10974097627.94334,10941577670.81993,10911063070.295,10879653592.23769,10849212601.51167,10817878983.52477,10787047730.29251,10756717154.06502,10725497325.48286,10695240359.11165,10664096390.5054,10775314007.31718,10744380616.61077,10710930909.92685,10679547283.9574,10647246726.16237,10615945472.75072,10583729805.66717,10552034069.37686,10520856375.28743,10488768044.26233,10457672722.89077,10425669282.55267,10559437536.99177,10527019635.1973,10493241632.89986,10461023029.27361,10427866480.68693,10395739012.96028,10362676385.84989,10330150734.97732,10298159971.00026,10265238226.36765,10233338598.16696,10200510775.27247,10330729118.83997,10299146192.88835,10265097889.23834,11928539254.32825,11035624852.2758,10776419395.07345,10545930715.66639,10322473034.62684,709916834275.6194)

# New "EBA_Adverse_2023_Euro Area" row (row 3) values
$row3Values = @(43437974611.26691,43089214286.76186,42819464544.53368,39302232384.96671,36384538308.79318,36148769562.17599,35909844294.94889,35676587278.28384,35440203044.94714,35534937581.16471,35325830696.86799,35123353901.93518,34928723832.51679,34731746623.20664,34538937923.35239,34343800348.94737,34151190732.31499,33961086191.15508,33768679949.07453,33580359744.60448,33389755806.77682,33568445001.21838,33393413691.6041,33211315872.4902,33035990996.91408,32858631950.89892,32684870400.62279,32509089315.11255,32335488221.02354,32164048519.80388,31990611061.22994,31820703965.39169,31648813516.50626,31899027553.33814,31728024395.38707,31550264272.50764,31379264770.38698,31206235100.19546,31036795489.93167,30865340815.27256,30696061390.39668,30528938163.37262,30359822344.10468,30194227869.0374,30026655663.15125,30335973907.93168,30167319778.5042,29995430881.90837,29828658821.35775,29659864843.1464,29494647565.85694,29327423888.9861,29162367361.43184,28999458533.48527,28834566394.46665,28673180852.92392,28509827267.5046,28801677191.81573,28641922454.81038,28475946890.25993,33007792496.98387,30460728255.69936,29670985810.24386,28963866837.49675,28279356492.252,2113570308515.559)

# Update row 2 ("Base") with the new values
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $row2Values[$i]
}

# Add new row 3 label
$ws.Range("A3").Value = "EBA_Adverse_2023_Euro Area"

# Add new row 3 values
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $row3Values[$i]
}
